# Modificación toneladas por día a toneladas totales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the daily-rate values in B9:B14 with formulas that multiply
# the original "tonnes per day" figure by the number of days in that
# month, turning them into "total tonnes" for the month.
$ws.Range("B9").Formula  = "=74*31"
$ws.Range("B10").Formula = "=53*28"
$ws.Range("B11").Formula = "=75*31"
$ws.Range("B12").Formula = "=139*30"
$ws.Range("B13").Formula = "=203*31"
$ws.Range("B14").Formula = "=164*30"

# Restore the active cell selection as recorded in the saved workbook.
$ws.Range("C10").Select()
